$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 / J1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy formatting (style) from existing header cell H1 so I1:J1 match s="1"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows 2..25 for columns I (9) and J (10)
$data = @(
    @(9,9),
    @(9,9),
    @(7,7),
    @(9,9),
    @(8,8),
    @(7,7),
    @(8,8),
    @(9,9),
    @(8,8),
    @(7,7),
    @(6,6),
    @(9,9),
    @(6,6),
    @(6,7),
    @(8,8),
    @(6,7),
    @(7,7),
    @(8,8),
    @(5,6),
    @(7,7),
    @(5,6),
    @(7,7),
    @(7,7),
    @(4,4)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
